$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D2:E51 so numeric-looking strings (e.g. "1.001") are not
# auto-converted to numbers by Excel, matching the original inlineStr text cells.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.423.80'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = '1.850.32'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '240.86'
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("D6").Value = '0.6304'
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '0.07703'
$ws.Range("E8").Value = '  +2.07%  '
$ws.Range("E9").Value = '  -0.34%  '
$ws.Range("D10").Value = '24.52'
$ws.Range("E10").Value = '  +0.31%  '
$ws.Range("E11").Value = '  +0.59%  '
$ws.Range("D12").Value = '1.854.39'
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").Value = '5.031'
$ws.Range("E13").Value = '  +0.72%  '
$ws.Range("D14").Value = '0.00001087'
$ws.Range("E14").Value = '  +8.08%  '
$ws.Range("D15").Value = '0.6803'
$ws.Range("E15").Value = '  -0.38%  '
$ws.Range("D16").Value = '83.77'
$ws.Range("E16").Value = '  +1.07%  '
$ws.Range("D17").Value = '2.107.36'
$ws.Range("D18").Value = '6.153'
$ws.Range("E18").Value = '  +0.32%  '
$ws.Range("D19").Value = '29.445.46'
$ws.Range("E19").Value = '  +0.16%  '
$ws.Range("D20").Value = '229.37'
$ws.Range("E20").Value = '  +0.84%  '
$ws.Range("E21").Value = '  +0.25%  '
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").Value = '7.452'
$ws.Range("E23").Value = '  -1.22%  '
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("D25").Value = '157.32'
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("E26").Value = '  -0.46%  '
$ws.Range("D27").Value = '8.378'
$ws.Range("E28").Value = '  +0.27%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '1.315'
$ws.Range("E29").Value = '  +4.41%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '1.468'
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("D31").Value = '0.05747'
$ws.Range("E31").Value = '  +1.25%  '
$ws.Range("D32").Value = '4.117'
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("D33").Value = '4.052'
$ws.Range("E33").Value = '  +0.89%  '
$ws.Range("D34").Value = '1.850'
$ws.Range("E34").Value = '  +0.35%  '
$ws.Range("E35").Value = '  +0.48%  '
$ws.Range("D36").Value = '0.7096'
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("D37").Value = '2.589'
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("D38").Value = '2.778'
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.01798'
$ws.Range("E39").Value = '  -0.88%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '1.226.09'
$ws.Range("E40").Value = '  -2.96%  '
$ws.Range("D41").Value = '6.469'
$ws.Range("E41").Value = '  +4.35%  '
$ws.Range("D42").Value = '0.9102'
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("D44").Value = '2.016.02'
$ws.Range("E44").Value = '  -0.73%  '
$ws.Range("D45").Value = '101.86'
$ws.Range("E45").Value = '  +0.46%  '
$ws.Range("D46").Value = '66.33'
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("D47").Value = '0.00000000120'
$ws.Range("E47").Value = '  +1.92%  '
$ws.Range("D48").Value = '7.161'
$ws.Range("E48").Value = '  +1.32%  '
$ws.Range("D49").Value = '0.4026'
$ws.Range("D50").Value = '9.034'
$ws.Range("E50").Value = '  -0.77%  '
$ws.Range("E51").Value = '  +0.45%  '

# Restore default style on the price/volume range (keeps values as text,
# removes the temporary text number-format override).
$priceVolRange.Style = "Normal"
